{"js": "// Remove the two outdated \"Current version limitations\" bullet points\n// (\"- Only deals with symmetric laminates\" and\n//  \"- No option to enforced balanced laminates\") from the Requirements /\n// limitations section, and relocate the \"_GoBack\" bookmark (which Word\n// stamps at the location of the most recent edit) from its old spot near\n// the end of the document onto the remaining\n// \"- Ply thickness remains constant\" bullet, right where the edit\n// actually happened.\n\nconst body = context.document.body;\n\n// Locate the two paragraphs to remove via their distinctive text.\nconst symResults = body.search(\"Only deals with symmetric laminates\", { matchCase: false });\nsymResults.load(\"items\");\nconst balResults = body.search(\"No option to enforced balanced laminates\", { matchCase: false });\nbalResults.load(\"items\");\n// The paragraph that survives and becomes the new \"_GoBack\" anchor.\nconst plyResults = body.search(\"Ply thickness remains constant\", { matchCase: false });\nplyResults.load(\"items\");\n\nawait context.sync();\n\nif (symResults.items.length === 0 || balResults.items.length === 0 || plyResults.items.length === 0) {\n  throw new Error(\"Could not locate the limitations bullet paragraphs to edit.\");\n}\n\n// Delete the whole paragraphs (not just the text) for both bullets.\nconst symPara = symResults.items[0].paragraphs.getFirst();\nconst balPara = balResults.items[0].paragraphs.getFirst();\nsymPara.delete();\nbalPara.delete();\n\n// Drop the old \"_GoBack\" bookmark wherever it currently lives (Word only\n// ever keeps a single \"_GoBack\" bookmark, marking the last edit spot).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Re-create it at the start of the remaining \"Ply thickness remains\n// constant\" bullet, matching where the real edit now lives.\nconst plyPara = plyResults.items[0].paragraphs.getFirst();\nconst plyStart = plyPara.getRange(\"Start\");\nplyStart.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Remove the two outdated \"Current version limitations\" bullet points\n# (\"- Only deals with symmetric laminates\" and\n#  \"- No option to enforced balanced laminates\") from the Requirements /\n# limitations section, and relocate the \"_GoBack\" bookmark (which Word\n# stamps at the location of the most recent edit) from its old spot near\n# the end of the document onto the remaining\n# \"- Ply thickness remains constant\" bullet, right where the edit\n# actually happened.\n\n$d = $word.ActiveDocument\n\n# --- Delete \"- Only deals with symmetric laminates\" paragraph ---\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Only deals with symmetric laminates\"\n$found = $find.Execute()\nif ($found) {\n    $para = $rng.Paragraphs(1)\n    $para.Range.Delete()\n}\n\n# --- Delete \"- No option to enforced balanced laminates\" paragraph ---\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.Text = \"No option to enforced balanced laminates\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $para2 = $rng2.Paragraphs(1)\n    $para2.Range.Delete()\n}\n\n# --- Move the \"_GoBack\" bookmark onto the surviving bullet ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$rng3 = $d.Content\n$find3 = $rng3.Find\n$find3.Text = \"Ply thickness remains constant\"\n$found3 = $find3.Execute()\nif ($found3) {\n    $plyPara = $rng3.Paragraphs(1)\n    $startPoint = $plyPara.Range.Duplicate\n    $startPoint.Collapse(1)\n    $d.Bookmarks.Add(\"_GoBack\", $startPoint)\n}\n"}
